$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.236.78'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '1.688.03'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.90'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.521'
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.97'
$ws.Range('E8').Value = '  +13.03%  '
$ws.Range('E9').Value = '  +3.89%  '
$ws.Range('E10').Value = '  +1.23%  '
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').Value = '1.926.60'
$ws.Range('E12').Value = '  +0.70%  '
$ws.Range('D13').Value = '1.691.57'
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('E14').Value = '  +2.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.555'
$ws.Range('E15').Value = '  +5.02%  '
$ws.Range('E16').Value = '  +2.19%  '
$ws.Range('D17').Value = '27.237.36'
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '239.90'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.17'
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('E20').Value = '  +1.75%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('E22').Value = '  +3.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.65'
$ws.Range('E23').Value = '  +4.83%  '
$ws.Range('E24').Value = '  -2.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.29'
$ws.Range('E25').Value = '  +1.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.32'
$ws.Range('E26').Value = '  +1.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.55'
$ws.Range('E27').Value = '  +2.43%  '
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0503'
$ws.Range('E30').Value = '  +0.90%  '
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.42'
$ws.Range('E32').Value = '  +2.85%  '
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.579.44'
$ws.Range('E33').Value = '  +6.73%  '
$ws.Range('E34').Value = '  +2.67%  '
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.953'
$ws.Range('E36').Value = '  +5.39%  '
$ws.Range('E37').Value = '  +3.35%  '
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('E39').Value = '  +0.47%  '
$ws.Range('E40').Value = '  +3.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.89'
$ws.Range('E41').Value = '  +3.48%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.65'
$ws.Range('E43').Value = '  -3.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.26'
$ws.Range('E44').Value = '  -2.44%  '
$ws.Range('D45').Value = '1.835.18'
$ws.Range('E45').Value = '  +0.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.788'
$ws.Range('E46').Value = '  +0.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '91.28'
$ws.Range('E47').Value = '  +0.83%  '
$ws.Range('E48').Value = '  +5.41%  '
$ws.Range('D49').Value = '0.0₆0108'
$ws.Range('E49').Value = '  +1.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.106'
$ws.Range('E50').Value = '  +3.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.21'
$ws.Range('E51').Value = '  +6.06%  '
